# Actualizacion automatica del tracker
# Appends 5 new result rows (55-59) to the tracker worksheet, mirroring the
# columns: event_id | fecha | jugador_A | jugador_B | pronostico | cuota | resultado | profit
# resultado / profit are still blank (pending) for these freshly-added matches,
# but stored as empty text (not a cleared cell) just like the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = 14687082; B = "2025-09-22"; C = "Alexander Bublik"; D = "Yibing Wu";          E = "Gana Alexander Bublik"; F = 1.4  },
    @{ A = 14733447; B = "2025-09-22"; C = "Blaise Bicknell";  D = "Patrick Maloney";     E = "Gana Blaise Bicknell";  F = 2.1  },
    @{ A = 14733450; B = "2025-09-22"; C = "Mats Rosenkranz";  D = "Cooper Williams";     E = "Gana Cooper Williams";  F = 1.83 },
    @{ A = 14738287; B = "2025-09-22"; C = "Remy Bertola";     D = "Kenny De Schepper";   E = "Gana Kenny De Schepper"; F = 2.62 },
    @{ A = 14738813; B = "2025-09-22"; C = "Leo Raquillet";    D = "Jelle Sels";          E = "Gana Leo Raquillet";    F = 2.62 }
)

$startRow = 55
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    # Leading apostrophe forces literal text so the ISO date string isn't
    # auto-converted into a date serial number by the COM value setter.
    $ws.Cells.Item($r, 2).Value = "'" + $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    # resultado / profit: not decided yet -> empty string (not a blank cell).
    # A lone apostrophe enters an empty literal-text value.
    $ws.Cells.Item($r, 7).Value = "'"
    $ws.Cells.Item($r, 8).Value = "'"
}
